$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 113
$wsExhibit.Range("F6").Value = 1889
$wsExhibit.Range("F8").Value = 60

# Sheet "全部类型" (All types) - fourth sheet, aggregated view shifted by one row
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 113
$wsAll.Range("F7").Value = 1889
$wsAll.Range("F9").Value = 60
